# Update cryptos list with latest scraped values (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.020.59'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.48%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.246.18'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.05%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.48'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.72%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '185.83'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.38%  '

$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.597'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.22%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.130'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.10%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.64'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.68%  '

$ws.Range("E11").Value = '  +1.14%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '3.813.87'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.16%  '

$ws.Range("E13").Value = '  +0.08%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.03'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.79%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '68.042.70'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.48%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000170'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.19%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.213.98'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.95%  '

$ws.Range("E18").Value = '  -0.03%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.52'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.19%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '394.60'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.44%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.61'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.03%  '

# Row 22/23: Dai and Litecoin swapped order
$ws.Range("B22").Value = 'Litecoin'
$ws.Range("C22").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '71.55'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.49%  '

$ws.Range("B23").Value = 'Dai'
$ws.Range("C23").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.23%  '

$ws.Range("E24").Value = '  +1.25%  '

$ws.Range("E25").Value = '  -0.01%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.188'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.30%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.68'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.86%  '

$ws.Range("E28").Value = '  +0.02%  '

$ws.Range("E29").Value = '  -0.02%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.66'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.21%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.89'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.39%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.09'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.17%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.26'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.51%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.999'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.12%  '

# Row 35/36: Monero and ImmutableX swapped order
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.50'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.85%  '

$ws.Range("B36").Value = 'Monero'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '162.04'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.16%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.91'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.58%  '

$ws.Range("E38").Value = '  -3.14%  '

# Row 39/40: EnergySwap and Filecoin swapped order
$ws.Range("B39").Value = 'Filecoin'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.64'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.87%  '

$ws.Range("B40").Value = 'EnergySwap'
$ws.Range("C40").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '26.44'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.35%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.54'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.81%  '

$ws.Range("E42").Value = '  -3.57%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0688'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.90%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.06'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.22%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '25.25'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.17%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.612.80'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.40%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '335.72'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.18%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0281'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.51%  '

$ws.Range("E49").Value = '  +3.61%  '

$ws.Range("E50").Value = '  -0.47%  '

# Row 51: Arweave replaced with ONDO
$ws.Range("B51").Value = 'ONDO'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.984'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.14%  '
